$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet: "Property1" -> "DataNode"
$ws.Name = "DataNode"

# Move the active selection to D37 (matches the saved view state in the diff)
$null = $ws.Range("D37").Select()
